$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(46006.01041666666, 0.37),
    @(46006.02083333334, 0.0),
    @(46006.03125, 0.0),
    @(46006.04166666666, 0.0),
    @(46006.05208333334, 0.45),
    @(46006.0625, 0.0),
    @(46006.07291666666, 0.0),
    @(46006.08333333334, 0.0),
    @(46006.09375, 0.0),
    @(46006.10416666666, 0.0),
    @(46006.11458333334, 0.0),
    @(46006.125, 0.0),
    @(46006.13541666666, 0.53),
    @(46006.14583333334, 0.0),
    @(46006.15625, 0.0),
    @(46006.16666666666, 0.0),
    @(46006.17708333334, 0.35),
    @(46006.1875, 0.0),
    @(46006.19791666666, 0.0),
    @(46006.20833333334, 0.0),
    @(46006.21875, 0.318),
    @(46006.22916666666, 0.356),
    @(46006.23958333334, 0.378),
    @(46006.25, 0.391),
    @(46006.26041666666, 2.766),
    @(46006.27083333334, 2.896),
    @(46006.28125, 3.435),
    @(46006.29166666666, 6.458),
    @(46006.30208333334, 50.951),
    @(46006.3125, 77.913),
    @(46006.32291666666, 113.179),
    @(46006.33333333334, 147.815),
    @(46006.34375, 309.069),
    @(46006.35416666666, 370.271),
    @(46006.36458333334, 431.587),
    @(46006.375, 491.184),
    @(46006.38541666666, 657.13),
    @(46006.39583333334, 720.688),
    @(46006.40625, 771.828),
    @(46006.41666666666, 818.959),
    @(46006.42708333334, 931.626),
    @(46006.4375, 968.916),
    @(46006.44791666666, 999.826),
    @(46006.45833333334, 1029.394),
    @(46006.46875, 1074.354),
    @(46006.47916666666, 1082.803),
    @(46006.48958333334, 1085.111),
    @(46006.5, 1077.492),
    @(46006.51041666666, 1048.745),
    @(46006.52083333334, 1020.233),
    @(46006.53125, 973.93),
    @(46006.54166666666, 935.342),
    @(46006.55208333334, 818.275),
    @(46006.5625, 761.353),
    @(46006.57291666666, 701.042),
    @(46006.58333333334, 635.541),
    @(46006.59375, 434.322),
    @(46006.60416666666, 363.956),
    @(46006.61458333334, 297.675),
    @(46006.625, 231.55),
    @(46006.63541666666, 87.756),
    @(46006.64583333334, 54.024),
    @(46006.65625, 36.094),
    @(46006.66666666666, 28.104),
    @(46006.67708333334, 10.964),
    @(46006.6875, 10.93),
    @(46006.69791666666, 11.478),
    @(46006.70833333334, 11.653),
    @(46006.71875, 3.86),
    @(46006.72916666666, 4.86),
    @(46006.73958333334, 0.0),
    @(46006.75, 2.86),
    @(46006.76041666666, 0.86),
    @(46006.77083333334, 0.0),
    @(46006.78125, 0.47),
    @(46006.79166666666, 0.0),
    @(46006.80208333334, 0.55),
    @(46006.8125, 0.0),
    @(46006.82291666666, 0.0),
    @(46006.83333333334, 0.0),
    @(46006.84375, 0.0),
    @(46006.85416666666, 0.0),
    @(46006.86458333334, 0.0),
    @(46006.875, 0.0),
    @(46006.88541666666, 0.53),
    @(46006.89583333334, 0.0),
    @(46006.90625, 0.0),
    @(46006.91666666666, 0.0),
    @(46006.92708333334, 0.45),
    @(46006.9375, 0.0),
    @(46006.94791666666, 0.0),
    @(46006.95833333334, 0.0),
    @(46006.96875, 0.0),
    @(46006.97916666666, 0.0),
    @(46006.98958333334, 0.0),
    @(46007.0, 0.0),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}